$wb = $excel.ActiveWorkbook

$rubric = $wb.Worksheets.Item("Rubric")
$grade = $wb.Worksheets.Item("Grade")

# Remove the "Indentation" note cell and revise the grade in column C, row 7
# on the Grade sheet (student grading removed from page 2).
$grade.Range("D7").ClearContents()
$grade.Range("C7").Value = 3

# Clear the (no longer needed) custom formatting on the sub-total rows.
$rubric.Range("A18:B18").ClearFormats()
$grade.Range("A18:C18").ClearFormats()

# Select single active cells on each sheet to match final saved selection state.
$rubric.Activate()
$rubric.Range("F11").Select()

# Make "Grade" the active (selected/visible) sheet, with its own selection.
$grade.Activate()
$grade.Range("E8").Select()

$wb.Save()
